$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3932893430561388
$ws.Range("C2").Value = 0.392707584027212

$ws.Range("B3").Value = 46.49187048396569
$ws.Range("C3").Value = 46.58285920086178

$ws.Range("B4").Value = 749.9267663892135
$ws.Range("C4").Value = 779.6498010712982

$ws.Range("B5").Value = 69.30000297728797
$ws.Range("C5").Value = 73.20397039299054

$ws.Range("B6").Value = 23825.80346679897
$ws.Range("C6").Value = 26321.73999868403

$ws.Range("B7").Value = 800.55176196153
$ws.Range("C7").Value = 1629.884265411335

$ws.Range("B8").Value = -3323.033283483964
$ws.Range("C8").Value = -134.8684947084942

$ws.Range("B9").Value = 573.5950408142635
$ws.Range("C9").Value = 585.9698573663768

$ws.Range("B10").Value = 3233.801199996195
$ws.Range("C10").Value = 5789.876948055387

$ws.Range("B11").Value = -1200.982765475774
$ws.Range("C11").Value = -377.4169906115071

$ws.Range("B12").Value = -8.476519280220081
$ws.Range("C12").Value = -8.447478717021379

$ws.Range("B13").Value = -4.084752035566997
$ws.Range("C13").Value = -4.339545630489408

$ws.Range("B14").Value = -4.025304721551356
$ws.Range("C14").Value = -4.28022145767423

$ws.Range("B15").Value = -0.7533606021254471
$ws.Range("C15").Value = -1.019443756518071
